# Auto-generated Excel COM-interop script to apply the Asura_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2322.862
$ws.Range("J112").Value = 2344.058
$ws.Range("L112").Value = 7032.174
$ws.Range("N112").Value = -9248.173999999999
$ws.Range("H134").Value = 100383.08
$ws.Range("J134").Value = 100383.08
$ws.Range("L134").Value = 100383.08
$ws.Range("N134").Value = -110523.08
$ws.Range("H137").Value = 1420.9429
$ws.Range("I137").Value = 1268.05
$ws.Range("J137").Value = 1624.8
$ws.Range("K137").Value = 3804.15
$ws.Range("L137").Value = 4874.4
$ws.Range("M137").Value = -1254.15
$ws.Range("N137").Value = -9974.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1139.6046
$ws.Range("I61").Value = 997.3871
$ws.Range("J61").Value = 1507
$ws.Range("K61").Value = 997.3871
$ws.Range("L61").Value = 1507
$ws.Range("M61").Value = -785.3871
$ws.Range("N61").Value = -1931
$ws.Range("H74").Value = 889.80554
$ws.Range("I74").Value = 845.7742
$ws.Range("K74").Value = 845.7742
$ws.Range("M74").Value = 28.22580000000005
$ws.Range("H77").Value = 889.80554
$ws.Range("I77").Value = 845.7742
$ws.Range("K77").Value = 4228.871
$ws.Range("M77").Value = 139.1289999999999
$ws.Range("H123").Value = 27972.75
$ws.Range("J123").Value = 27972.75
$ws.Range("L123").Value = 27972.75
$ws.Range("N123").Value = -37772.75
$ws.Range("H132").Value = 2087.1135
$ws.Range("I132").Value = 1337.15
$ws.Range("K132").Value = 4011.45
$ws.Range("M132").Value = -1481.45
$ws.Range("H136").Value = 1139.6046
$ws.Range("I136").Value = 997.3871
$ws.Range("J136").Value = 1507
$ws.Range("K136").Value = 2992.1613
$ws.Range("L136").Value = 4521
$ws.Range("M136").Value = -442.1613000000002
$ws.Range("N136").Value = -9621

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 30148.334
$ws.Range("I107").Value = 37875.145
$ws.Range("J107").Value = 3104.5
$ws.Range("K107").Value = 37875.145
$ws.Range("L107").Value = 3104.5
$ws.Range("M107").Value = -35955.145
$ws.Range("N107").Value = -6944.5
$ws.Range("H134").Value = 1902.0652
$ws.Range("I134").Value = 1601.7632
$ws.Range("J134").Value = 3328.5
$ws.Range("K134").Value = 4805.2896
$ws.Range("L134").Value = 9985.5
$ws.Range("M134").Value = -2270.2896
$ws.Range("N134").Value = -15055.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2377.342
$ws.Range("I31").Value = 1485.28
$ws.Range("J31").Value = 4092.8462
$ws.Range("K31").Value = 1485.28
$ws.Range("L31").Value = 4092.8462
$ws.Range("M31").Value = -1190.28
$ws.Range("N31").Value = -4682.8462
$ws.Range("H34").Value = 2377.342
$ws.Range("I34").Value = 1485.28
$ws.Range("J34").Value = 4092.8462
$ws.Range("K34").Value = 1485.28
$ws.Range("L34").Value = 4092.8462
$ws.Range("M34").Value = -1283.28
$ws.Range("N34").Value = -4496.8462
$ws.Range("H58").Value = 662620.7
$ws.Range("I58").Value = 975730.7
$ws.Range("J58").Value = 1610.5555
$ws.Range("K58").Value = 975730.7
$ws.Range("L58").Value = 1610.5555
$ws.Range("M58").Value = -975527.7
$ws.Range("N58").Value = -2016.5555
$ws.Range("H132").Value = 357077.5
$ws.Range("I132").Value = 541753.8
$ws.Range("J132").Value = 1930.6923
$ws.Range("K132").Value = 1625261.4
$ws.Range("L132").Value = 5792.0769
$ws.Range("M132").Value = -1622731.4
$ws.Range("N132").Value = -10852.0769
$ws.Range("H134").Value = 1359.8462
$ws.Range("I134").Value = 1147.6428
$ws.Range("K134").Value = 3442.9284
$ws.Range("M134").Value = -907.9284000000002
$ws.Range("H136").Value = 662620.7
$ws.Range("I136").Value = 975730.7
$ws.Range("J136").Value = 1610.5555
$ws.Range("K136").Value = 2927192.1
$ws.Range("L136").Value = 4831.666499999999
$ws.Range("M136").Value = -2924642.1
$ws.Range("N136").Value = -9931.666499999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 189.28572
$ws.Range("I4").Value = 189.28572
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 567.85716
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -455.85716
$ws.Range("N4").Value = ""
$ws.Range("H11").Value = 275
$ws.Range("I11").Value = 85.28570999999999
$ws.Range("K11").Value = 255.85713
$ws.Range("M11").Value = -115.85713
$ws.Range("H123").Value = 4266.6665
$ws.Range("I123").Value = 800
$ws.Range("J123").Value = 6000
$ws.Range("K123").Value = 2400
$ws.Range("L123").Value = 18000
$ws.Range("M123").Value = 50
$ws.Range("N123").Value = -22900
$ws.Range("H129").Value = 3573007
$ws.Range("I129").Value = 899.75
$ws.Range("J129").Value = 5001850
$ws.Range("K129").Value = 2699.25
$ws.Range("L129").Value = 15005550
$ws.Range("M129").Value = 2300.75
$ws.Range("N129").Value = -15015550
$ws.Range("H130").Value = 3000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 9000
$ws.Range("M130").Value = ""
$ws.Range("N130").Value = -19040
$ws.Range("H131").Value = 14098783
$ws.Range("I131").Value = 1526.25
$ws.Range("J131").Value = 15888911
$ws.Range("K131").Value = 4578.75
$ws.Range("L131").Value = 47666733
$ws.Range("M131").Value = 461.25
$ws.Range("N131").Value = -47676813
$ws.Range("H133").Value = 4957.857
$ws.Range("I133").Value = 2015
$ws.Range("J133").Value = 5448.3335
$ws.Range("K133").Value = 6045
$ws.Range("L133").Value = 16345.0005
$ws.Range("M133").Value = -985
$ws.Range("N133").Value = -26465.0005
$ws.Range("H134").Value = 4316.921
$ws.Range("I134").Value = 1774.6154
$ws.Range("J134").Value = 5638.92
$ws.Range("K134").Value = 5323.8462
$ws.Range("L134").Value = 16916.76
$ws.Range("M134").Value = -253.8462
$ws.Range("N134").Value = -27056.76
$ws.Range("H136").Value = 3629.889
$ws.Range("I136").Value = 1086.25
$ws.Range("J136").Value = 7329.727
$ws.Range("K136").Value = 3258.75
$ws.Range("L136").Value = 21989.181
$ws.Range("M136").Value = 1841.25
$ws.Range("N136").Value = -32189.181
$ws.Range("H137").Value = 16670240
$ws.Range("I137").Value = 3843.3333
$ws.Range("J137").Value = 19611368
$ws.Range("K137").Value = 11529.9999
$ws.Range("L137").Value = 58834104
$ws.Range("M137").Value = -6429.999899999999
$ws.Range("N137").Value = -58844304
$ws.Range("H138").Value = 1987.9
$ws.Range("I138").Value = 989.8333
$ws.Range("J138").Value = 3485
$ws.Range("K138").Value = 2969.4999
$ws.Range("L138").Value = 10455
$ws.Range("M138").Value = 2170.5001
$ws.Range("N138").Value = -20735
$ws.Range("H139").Value = 2342.7742
$ws.Range("I139").Value = 2093.6365
$ws.Range("J139").Value = 2951.7778
$ws.Range("K139").Value = 6280.9095
$ws.Range("L139").Value = 8855.3334
$ws.Range("M139").Value = -1140.9095
$ws.Range("N139").Value = -19135.3334
$ws.Range("H140").Value = 1768.4482
$ws.Range("I140").Value = 941.7727
$ws.Range("J140").Value = 4366.5713
$ws.Range("K140").Value = 2825.3181
$ws.Range("L140").Value = 13099.7139
$ws.Range("M140").Value = 2354.6819
$ws.Range("N140").Value = -23459.7139
$ws.Range("H141").Value = 4966.647
$ws.Range("I141").Value = 3725.6155
$ws.Range("J141").Value = 9000
$ws.Range("K141").Value = 11176.8465
$ws.Range("L141").Value = 27000
$ws.Range("M141").Value = -5996.8465
$ws.Range("N141").Value = -37360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 26230.46
$ws.Range("J5").Value = 26230.46
$ws.Range("L5").Value = 26230.46
$ws.Range("N5").Value = -26454.46
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H108").Value = 34000
$ws.Range("J108").Value = 34000
$ws.Range("L108").Value = 34000
$ws.Range("N108").Value = -41680
$ws.Range("H132").Value = 1360.4193
$ws.Range("I132").Value = 794.3182
$ws.Range("J132").Value = 2744.2222
$ws.Range("K132").Value = 2382.9546
$ws.Range("L132").Value = 8232.6666
$ws.Range("M132").Value = 147.0454
$ws.Range("N132").Value = -13292.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9926.143
$ws.Range("J2").Value = 9926.143
$ws.Range("L2").Value = 9926.143
$ws.Range("N2").Value = -10150.143
$ws.Range("H100").Value = 15425
$ws.Range("I100").Value = 50000
$ws.Range("J100").Value = 3900
$ws.Range("K100").Value = 50000
$ws.Range("L100").Value = 3900
$ws.Range("M100").Value = -49459
$ws.Range("N100").Value = -4982
$ws.Range("H132").Value = 2475.0532
$ws.Range("I132").Value = 2043.017
$ws.Range("J132").Value = 4068.1875
$ws.Range("K132").Value = 6129.051
$ws.Range("L132").Value = 12204.5625
$ws.Range("M132").Value = -3599.051
$ws.Range("N132").Value = -17264.5625
$ws.Range("H136").Value = 2486.4045
$ws.Range("I136").Value = 2589.3623
$ws.Range("K136").Value = 7768.086899999999
$ws.Range("M136").Value = -5218.086899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22416.445
$ws.Range("J123").Value = 22416.445
$ws.Range("L123").Value = 22416.445
$ws.Range("N123").Value = -32216.445
$ws.Range("H132").Value = 926.73846
$ws.Range("I132").Value = 615.5282999999999
$ws.Range("K132").Value = 1846.5849
$ws.Range("M132").Value = 683.4151000000002
$ws.Range("H136").Value = 2046.3611
$ws.Range("I136").Value = 1643.7587
$ws.Range("K136").Value = 4931.2761
$ws.Range("M136").Value = -2381.2761

